$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray youtube-link cell that used to sit in column D of row 10.
# (Its shared-string entry disappears entirely once no cell references it.)
$ws.Range("D10").ClearContents()

# Row 12: rewrite the description text (typo fix "ajouté" -> "ajouter" plus a
# new trailing sentence about the bottle-removal movement).
$ws.Range("B12").Value2 = "Une des deux pages de mon application mobile est terminée, j'ai commencé à faire la deuxième page. Je peux déjà ajouter un mouvement quand j'ajoute des bouteilles dans ma cave. Elle reprend aussi le fait qu'une personne enlève une bouteille et cela crée un mouvement dans la BD."

# Add a new row 13 for the day's extra task, copying the formatting (date
# format on A, wrap-text on B) from the row above so no new cell styles are
# introduced.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value2 = 43237
$ws.Range("B13").Value2 = "J'ai fait de la documentation, j'ai ajouté mes maquettes et aimélioré des uses cases et scénarios "
$ws.Range("C13").Value2 = "2 heures"

# Update the view: zoom back to 100% and move the selection to the next
# empty row in column C.
$win = $ws.Application.ActiveWindow
$win.Zoom = 100
$ws.Range("C14").Select() | Out-Null
